$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (days since 1899-12-30).
# Update all data rows (2-348) from serial 45175 (2023-09-06) to 45177 (2023-09-08).
$ws.Range("C2:C348").Value = 45177
